$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (PRICE -> new col -> QUANTITY shifts right)
$ws.Columns("D").EntireColumn.Insert()

# New column takes the same width as the PRICE column next to it
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(3).ColumnWidth

# Header for the newly inserted column
$ws.Cells.Item(1, 4).Value = "TRIGGER_PRICE"

# Trigger price values for each order row
$ws.Cells.Item(2, 4).Value = 2.2
$ws.Cells.Item(3, 4).Value = 3.1
$ws.Cells.Item(4, 4).Value = 5.5
$ws.Cells.Item(5, 4).Value = 1.4

# All orders are now BUY orders (SELL no longer used)
$ws.Cells.Item(3, 6).Value = "BUY"
$ws.Cells.Item(4, 6).Value = "BUY"

# Update selection to mirror the author's final cursor position
$ws.Range("F6").Select()
